$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 240-241, shifting existing rows 240:339 down to 242:341
$ws.Range("A240:A241").EntireRow.Insert()

# Fill the two newly inserted rows with the new data records.
# Row 240: new "Zafiro rojo" record dated 44846
$ws.Range("A240").Value = 7
$ws.Range("B240").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C240").Value = "Ñuble"
$ws.Range("D240").Value = 44846
$ws.Range("E240").Value = 16
$ws.Range("F240").Value = 100112002
$ws.Range("G240").Value = "Pimiento"
$ws.Range("H240").Value = "Zafiro rojo"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 100
$ws.Range("K240").Value = 22000
$ws.Range("L240").Value = 23000
$ws.Range("M240").Value = 22500
$ws.Range("N240").Value = '$/caja 15 kilos'
$ws.Range("O240").Value = "Región de Arica y Parinacota"
$ws.Range("P240").Value = 1500
$ws.Range("Q240").Value = 15
$ws.Range("R240").Value = "Hortaliza"

# Row 241: new "Zafiro verde" record dated 44846
$ws.Range("A241").Value = 7
$ws.Range("B241").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C241").Value = "Ñuble"
$ws.Range("D241").Value = 44846
$ws.Range("E241").Value = 16
$ws.Range("F241").Value = 100112002
$ws.Range("G241").Value = "Pimiento"
$ws.Range("H241").Value = "Zafiro verde"
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 100
$ws.Range("K241").Value = 22000
$ws.Range("L241").Value = 23000
$ws.Range("M241").Value = 22500
$ws.Range("N241").Value = '$/caja 15 kilos'
$ws.Range("O241").Value = "Región de Arica y Parinacota"
$ws.Range("P241").Value = 1500
$ws.Range("Q241").Value = 15
$ws.Range("R241").Value = "Hortaliza"
